$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.133.78'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '3.139.21'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.512'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.77%  '
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('E11').Value = '  +6.20%  '
$ws.Range('E12').Value = '  +3.14%  '
$ws.Range('D13').Value = '3.678.09'
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('E15').Value = '  +4.96%  '
$ws.Range('D16').Value = '58.241.09'
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.91%  '
$ws.Range('D18').Value = '3.128.58'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.61%  '
$ws.Range('E20').Value = '  +4.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.99%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.31%  '
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +12.66%  '
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.91%  '
$ws.Range('E32').Value = '  +3.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.55%  '
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '161.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.74'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.72%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.644.36'
$ws.Range('E40').Value = '  +10.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.27'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0679'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.701'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('E45').Value = '  +4.96%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  +4.80%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.101'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.41%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.978'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.38'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.752'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.49%  '
